$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A4").Value = "A 13"
$ws.Range("B4").Value = "B 13"
$ws.Range("C4").Value = "C 13"
$ws.Range("D4").Value = "D 13"
